$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Content fix: remove the erroneous trailing `Cohort` column from the
#     CasesTab query stored in B2 (cohort_description was never requested
#     by this test case's WHERE clause / intent). ---
$newQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
MATCH (c)<--(diag:diagnosis)
MATCH (samp:sample)-->(c) 
  MATCH (f:file)-[*]->(c)
    WHERE f.file_format IN ["bam"]  
OPTIONAL MATCH (co:cohort)<-[*]-(c)
  WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(demo.weight, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment`
'@
$ws.Range("B2").Value2 = $newQuery

# --- View state: the query was re-reviewed/re-saved with the sheet
#     zoomed to 85% and the selection left on the edited cell (B2)
#     instead of the previous B4/55% view. ---
$excel.ActiveWindow.Zoom = 85
$ws.Range("B2").Select()
